# Insert two new data rows at row 391 (pushing the existing rows 391.. down by two,
# so old row 391 becomes row 393, ..., old row 479 becomes row 481).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A391:A392").EntireRow.Insert()

# New row 391
$ws.Cells.Item(391, 1).Value2 = 11
$ws.Cells.Item(391, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(391, 3).Value2 = "Bíobío"
$ws.Cells.Item(391, 4).Value2 = 44785
$ws.Cells.Item(391, 5).Value2 = 8
$ws.Cells.Item(391, 6).Value2 = 100112020
$ws.Cells.Item(391, 7).Value2 = "Tomate"
$ws.Cells.Item(391, 8).Value2 = "Larga vida"
$ws.Cells.Item(391, 9).Value2 = "Primera"
$ws.Cells.Item(391, 10).Value2 = 500
$ws.Cells.Item(391, 11).Value2 = 8500
$ws.Cells.Item(391, 12).Value2 = 9000
$ws.Cells.Item(391, 13).Value2 = 8700
$ws.Cells.Item(391, 14).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(391, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(391, 16).Value2 = 483
$ws.Cells.Item(391, 17).Value2 = 18
$ws.Cells.Item(391, 18).Value2 = "Hortaliza"

# New row 392
$ws.Cells.Item(392, 1).Value2 = 11
$ws.Cells.Item(392, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(392, 3).Value2 = "Bíobío"
$ws.Cells.Item(392, 4).Value2 = 44785
$ws.Cells.Item(392, 5).Value2 = 8
$ws.Cells.Item(392, 6).Value2 = 100112020
$ws.Cells.Item(392, 7).Value2 = "Tomate"
$ws.Cells.Item(392, 8).Value2 = "Larga vida"
$ws.Cells.Item(392, 9).Value2 = "Segunda"
$ws.Cells.Item(392, 10).Value2 = 300
$ws.Cells.Item(392, 11).Value2 = 8000
$ws.Cells.Item(392, 12).Value2 = 8000
$ws.Cells.Item(392, 13).Value2 = 8000
$ws.Cells.Item(392, 14).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(392, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(392, 16).Value2 = 444
$ws.Cells.Item(392, 17).Value2 = 18
$ws.Cells.Item(392, 18).Value2 = "Hortaliza"
